$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.460.26'
$ws.Range('E2').Value = '  +0.47%  '

$ws.Range('D3').Value = '2.269.63'
$ws.Range('E3').Value = '  -0.02%  '

$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').Value = '307.97'
$ws.Range('E5').Value = '  +0.99%  '

$ws.Range('D6').Value = '99.30'
$ws.Range('E6').Value = '  +2.12%  '

$ws.Range('D7').Value = '0.526'
$ws.Range('E7').Value = '  -0.60%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').Value = '0.494'
$ws.Range('E9').Value = '  +0.61%  '

$ws.Range('D10').Value = '35.59'
$ws.Range('E10').Value = '  -0.03%  '

$ws.Range('D11').Value = '0.0820'
$ws.Range('E11').Value = '  +3.07%  '

$ws.Range('E12').Value = '  +1.86%  '

$ws.Range('D13').Value = '6.81'
$ws.Range('E13').Value = '  +2.30%  '

$ws.Range('D14').Value = '2.619.10'
$ws.Range('E14').Value = '  -0.04%  '

$ws.Range('D15').Value = '14.65'
$ws.Range('E15').Value = '  +1.89%  '

$ws.Range('D16').Value = '2.262.23'
$ws.Range('E16').Value = '  -0.14%  '

$ws.Range('D17').Value = '0.787'
$ws.Range('E17').Value = '  -1.03%  '

$ws.Range('D18').Value = '42.264.31'
$ws.Range('E18').Value = '  +0.26%  '

$ws.Range('D19').Value = '12.31'
$ws.Range('E19').Value = '  -1.44%  '

$ws.Range('D20').Value = '0.0₃0909'
$ws.Range('E20').Value = '  +0.07%  '

$ws.Range('D21').Value = '5.99'
$ws.Range('E21').Value = '  +0.41%  '

$ws.Range('D22').Value = '67.38'
$ws.Range('E22').Value = '  -0.28%  '

$ws.Range('D23').Value = '237.40'
$ws.Range('E23').Value = '  -0.15%  '

$ws.Range('D24').Value = '2.59'
$ws.Range('E24').Value = '  +0.29%  '

$ws.Range('D25').Value = '1.96'
$ws.Range('E25').Value = '  +0.34%  '

$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('D27').Value = '38.41'
$ws.Range('E27').Value = '  +2.35%  '

$ws.Range('D28').Value = '23.58'
$ws.Range('E28').Value = '  -1.01%  '

$ws.Range('D29').Value = '2.15'
$ws.Range('E29').Value = '  +1.33%  '

$ws.Range('D30').Value = '9.57'
$ws.Range('E30').Value = '  +0.64%  '

$ws.Range('D31').Value = '167.26'
$ws.Range('E31').Value = '  +4.57%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.22'
$ws.Range('E32').Value = '  -0.72%  '

$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').Value = '0.997'
$ws.Range('E33').Value = '  -0.14%  '

$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '3.08'
$ws.Range('E34').Value = '  -2.28%  '

$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').Value = '17.63'
$ws.Range('E35').Value = '  +3.10%  '

$ws.Range('D36').Value = '0.0726'
$ws.Range('E36').Value = '  -2.09%  '

$ws.Range('E37').Value = '  +1.25%  '

$ws.Range('D38').Value = '0.115'
$ws.Range('E38').Value = '  +0.42%  '

$ws.Range('E39').Value = '  -1.62%  '

$ws.Range('D40').Value = '1.80'
$ws.Range('E40').Value = '  -1.30%  '

$ws.Range('D41').Value = '4.15'
$ws.Range('E41').Value = '  +2.29%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '19.01'
$ws.Range('E42').Value = '  +0.72%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.944.21'
$ws.Range('E43').Value = '  -2.50%  '

$ws.Range('D44').Value = '0.0283'
$ws.Range('E44').Value = '  -1.05%  '

$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = '2.22'
$ws.Range('E45').Value = '  -8.55%  '

$ws.Range('D46').Value = '2.93'
$ws.Range('E46').Value = '  -0.41%  '

$ws.Range('D47').Value = '9.76'
$ws.Range('E47').Value = '  -2.06%  '

$ws.Range('D48').Value = '54.27'
$ws.Range('E48').Value = '  +1.93%  '

$ws.Range('D49').Value = '2.484.59'
$ws.Range('E49').Value = '  -0.20%  '

$ws.Range('D50').Value = '71.74'
$ws.Range('E50').Value = '  -0.46%  '

$ws.Range('D51').Value = '91.88'
$ws.Range('E51').Value = '  +0.53%  '
